# Refresh "想去人数" (interest count, column F) figures across all sheets,
# matching the upstream data regeneration (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 771
$ws.Range("F3").Value  = 2795
$ws.Range("F4").Value  = 1330
$ws.Range("F6").Value  = 1961
$ws.Range("F9").Value  = 604
$ws.Range("F12").Value = 11641
$ws.Range("F13").Value = 6616
$ws.Range("F16").Value = 418
$ws.Range("F20").Value = 917
$ws.Range("F21").Value = 80
$ws.Range("F23").Value = 925
$ws.Range("F24").Value = 3647
$ws.Range("F28").Value = 167
$ws.Range("F31").Value = 22
$ws.Range("F32").Value = 268
$ws.Range("F33").Value = 302
$ws.Range("F34").Value = 5013
$ws.Range("F36").Value = 1239
$ws.Range("F37").Value = 231
$ws.Range("F38").Value = 442
$ws.Range("F39").Value = 198
$ws.Range("F40").Value = 538

# --- Sheet "演出" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 1
$ws.Range("F24").Value = 39

# --- Sheet "本地生活" -------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9046
$ws.Range("F3").Value = 499
$ws.Range("F4").Value = 1829

# --- Sheet "全部类型" -------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 499
$ws.Range("F3").Value  = 1829
$ws.Range("F4").Value  = 771
$ws.Range("F5").Value  = 2795
$ws.Range("F9").Value  = 1330
$ws.Range("F13").Value = 48
$ws.Range("F15").Value = 604
$ws.Range("F18").Value = 11641
$ws.Range("F20").Value = 6616
$ws.Range("F24").Value = 418
$ws.Range("F28").Value = 80
$ws.Range("F30").Value = 925
$ws.Range("F31").Value = 3647
$ws.Range("F34").Value = 167
$ws.Range("F36").Value = 22
$ws.Range("F37").Value = 268
$ws.Range("F41").Value = 1239
$ws.Range("F42").Value = 231
$ws.Range("F44").Value = 198
$ws.Range("F45").Value = 538
